$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# New test case row: Comptroller source object profile (Login_Logout / Logout)
# Duplicate row 16's formatting onto the new row 17, then fill in the values.
$ws.Rows("16").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Rows("17").RowHeight = 22

$ws.Range("A17").Value = "MDOT"
$ws.Range("B17").Value = "Login_Logout"
$ws.Range("C17").Value = "Logout"
$ws.Range("D17").Value = "Yes"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1

# Update the active view to the newly added row
$ws.Range("A17:F17").Select()
$excel.ActiveWindow.ScrollRow = 5
